$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new row at position 3, shifting the existing rows 3-5 down to 4-6.
    $ws.Rows.Item(3).Insert()

    # Restore the bordered index-column style (Insert's auto-copied style drops
    # the border) so the new A3 matches the other numbered rows.
    $ws.Range("A3").Borders.LineStyle = 1

    # Fill in the newly inserted row 3 with the new exhibition entry.
    # Force column B to remain plain text (matching the sibling date cells)
    # instead of letting Excel auto-convert the "2024-08-10" string to a date,
    # then drop back to the plain "Normal" style (no cell format override),
    # matching the unstyled text cells used by the rest of the column.
    $ws.Range("B3").NumberFormat = "@"
    $ws.Range("B3").Value = "2024-08-10"
    $ws.Range("B3").Style = "Normal"
    $ws.Range("C3").Value = "丽水·未来城次元同好会免费展"
    $ws.Range("D3").Value = "中广未来城1019号、1020号 丽水经开未来城"
    $ws.Range("E3").Value = "2024.08.10 12:30-08.11 17:00"
    $ws.Range("F3").Value = 1
    $ws.Range("G3").Value = 29.9
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=90282"
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202408/CW6Lw0Z11722583407396.jpeg"

    # Renumber the index column (A) sequentially for all data rows (2-6).
    $ws.Range("A2").Value = 1
    $ws.Range("A3").Value = 2
    $ws.Range("A4").Value = 3
    $ws.Range("A5").Value = 4
    $ws.Range("A6").Value = 5
}
